# Updates odds values for the two match rows (row 3 and row 6) on Sheet1,
# reflecting refreshed FlashScore odds data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 (CSKA Sofia - Botev Vratsa) ---
$ws.Range("G3").Value  = 1.25
$ws.Range("H3").Value  = 5.5
$ws.Range("I3").Value  = 11
$ws.Range("M3").Value  = 1.06
$ws.Range("N3").Value  = 10
$ws.Range("U3").Value  = 2.63
$ws.Range("V3").Value  = 1.44
$ws.Range("Y3").Value  = 10
$ws.Range("Z3").Value  = 7
$ws.Range("AD3").Value = 11
$ws.Range("AH3").Value = 21
$ws.Range("AJ3").Value = 34
$ws.Range("AM3").Value = 101
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 11

# --- Row 6 (Al Orubah - Al Shabab) ---
$ws.Range("G6").Value  = 3.5
$ws.Range("I6").Value  = 2.05
$ws.Range("L6").Value  = 2.75
$ws.Range("M6").Value  = 1.06
$ws.Range("N6").Value  = 8
$ws.Range("W6").Value  = 10
$ws.Range("X6").Value  = 17
$ws.Range("AA6").Value = 29
$ws.Range("AC6").Value = 8.5
$ws.Range("AI6").Value = 9.5
$ws.Range("AJ6").Value = 9.5
$ws.Range("AL6").Value = 19
$ws.Range("AR6").Value = 81

$wb.Save()
